$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.896318674087524
$ws.Range("B1").Value = 4.970314979553223
$ws.Range("C1").Value = 3.58087682723999
$ws.Range("D1").Value = 1.190950870513916
$ws.Range("E1").Value = 0.7839040756225586
